# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Reordena la columna "Periodo Mora" (E16:E20) de orden descendente a
#   ascendente: 1705,1704,1702,1701,1612 -> 1612,1701,1702,1704,1705
# - Actualiza la columna "Valor Mora" (G16:G20) de 689455 a 781242

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1612", "1701", "1702", "1704", "1705")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}

$ws.Range("G16:G20").Value = 781242
